$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 8 - S.NO 3, Dream uniforms, new PO
$ws.Cells.Item(8, 1).Value = 3
$ws.Cells.Item(8, 2).Value = "Dream uniforms"
$ws.Cells.Item(8, 3).Value = "lpo/Dream uniforms/DU/PO/001Test"
Set-TextValue 8 4 "4657-4756"
$ws.Cells.Item(8, 5).Value = "Shirt"
$ws.Cells.Item(8, 6).Value = "Dream_uniform_2"
Set-TextValue 8 7 "100"
Set-TextValue 8 8 "10"
$ws.Cells.Item(8, 9).Value = "abcd cut on 10/1/2019--19:12"
Set-TextValue 8 10 "0"
$ws.Cells.Item(8, 11).Value = " stitched on "
Set-TextValue 8 12 "0"
$ws.Cells.Item(8, 13).Value = " finalized on "
Set-TextValue 8 14 "0"
$ws.Cells.Item(8, 15).Value = " packing on "
Set-TextValue 8 16 "100"

# Row 9 - S.NO 4, Affan uniforms co.
$ws.Cells.Item(9, 1).Value = 4
$ws.Cells.Item(9, 2).Value = "Affan uniforms co."
$ws.Cells.Item(9, 3).Value = "lpo/Affan uniforms co./500"
Set-TextValue 9 4 "4757-4806"
$ws.Cells.Item(9, 5).Value = "Shirt"
$ws.Cells.Item(9, 6).Value = "Geo Anchor"
Set-TextValue 9 7 "50"
Set-TextValue 9 8 "0"
$ws.Cells.Item(9, 9).Value = " cut on "
Set-TextValue 9 10 "0"
$ws.Cells.Item(9, 11).Value = " stitched on "
Set-TextValue 9 12 "0"
$ws.Cells.Item(9, 13).Value = " finalized on "
Set-TextValue 9 14 "50"
$ws.Cells.Item(9, 15).Value = "ali packing on 14/1/2019--14:45"
Set-TextValue 9 16 "0"

# Row 10 - S.NO 5, Affan uniforms co.
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = "Affan uniforms co."
$ws.Cells.Item(10, 3).Value = "lpo/Affan uniforms co./500"
Set-TextValue 10 4 "4807-4856"
$ws.Cells.Item(10, 5).Value = "Shirt"
$ws.Cells.Item(10, 6).Value = "Geo tv anchor female"
Set-TextValue 10 7 "50"
Set-TextValue 10 8 "0"
$ws.Cells.Item(10, 9).Value = " cut on "
Set-TextValue 10 10 "0"
$ws.Cells.Item(10, 11).Value = " stitched on "
Set-TextValue 10 12 "0"
$ws.Cells.Item(10, 13).Value = " finalized on "
Set-TextValue 10 14 "17"
$ws.Cells.Item(10, 15).Value = "ahmed packing on 14/1/2019--14:46"
Set-TextValue 10 16 "33"

# Row 12 - totals
$ws.Cells.Item(12, 7).Value = "Total Cut"
$ws.Cells.Item(12, 8).Value = 60
$ws.Cells.Item(12, 9).Value = "Total Stitched"
$ws.Cells.Item(12, 10).Value = 22
$ws.Cells.Item(12, 11).Value = "Total Finished"
$ws.Cells.Item(12, 12).Value = 14
$ws.Cells.Item(12, 13).Value = "Total Packed"
$ws.Cells.Item(12, 14).Value = 81
$ws.Cells.Item(12, 15).Value = "Total Delivered"
$ws.Cells.Item(12, 16).Value = 81
